# Update `vendors` sample data
# Adds 7 sample vendor rows (name + vendor id) to the "initialize_vendors"
# sheet, sizes the columns to fit the new content, and makes that sheet
# the active/selected one (matching where the editor left the cursor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("initialize_vendors")

# --- sample vendor data -----------------------------------------------
$vendors = @(
    @{ Name = "ProQuest";           Id = 1 },
    @{ Name = "EBSCO";              Id = 2 },
    @{ Name = "Gale";               Id = 3 },
    @{ Name = "iG Publishing/BEP";  Id = 4 },
    @{ Name = "Ebook Library";      Id = 5 },
    @{ Name = "Ebrary";             Id = 6 },
    @{ Name = "MyiLibrary";         Id = 7 }
)

$row = 2
foreach ($vendor in $vendors) {
    $ws.Range("A$row").Value = $vendor.Name
    $ws.Range("C$row").Value = $vendor.Id
    $row++
}

# --- resize columns to fit the new data --------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.0
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 8.833333333333332

# --- make this sheet the active one, with D7 selected ------------------
$ws.Activate() | Out-Null
$ws.Range("D7").Select() | Out-Null
